# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F-column "time_taken" timestamps on the "data" sheet ---
$newTimes = @(
    "2021-10-05 14:35:53.730159",
    "2021-10-05 14:35:53.730167",
    "2021-10-05 14:35:53.730171",
    "2021-10-05 14:35:53.730173",
    "2021-10-05 14:35:53.730176",
    "2021-10-05 14:35:53.730179",
    "2021-10-05 14:35:53.730182",
    "2021-10-05 14:35:53.730184",
    "2021-10-05 14:35:53.730187",
    "2021-10-05 14:35:53.730190",
    "2021-10-05 14:35:53.730192",
    "2021-10-05 14:35:53.730195",
    "2021-10-05 14:35:53.730198",
    "2021-10-05 14:35:53.730200",
    "2021-10-05 14:35:53.730203",
    "2021-10-05 14:35:53.730205",
    "2021-10-05 14:35:53.730208",
    "2021-10-05 14:35:53.730211",
    "2021-10-05 14:35:53.730213",
    "2021-10-05 14:35:53.730216",
    "2021-10-05 14:35:53.730218",
    "2021-10-05 14:35:53.730221",
    "2021-10-05 14:35:53.730223",
    "2021-10-05 14:35:53.730226"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" worksheet, positioned after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Vascular Malformations_Somatic"
$metaSheet.Range("C2").Value = 3181
# Force the version string to stay text (it would otherwise be parsed as
# the number 1.6) by using the leading-apostrophe text-entry convention.
$metaSheet.Range("D2").Value = "'1.6"
$metaSheet.Range("E2").Value = "2021-09-02T05:37:57.210299Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:53.726450"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3181/?format=json"

# Apply the same header/bordered style used on the "data" sheet's header row
# (and its A2 index cell) to the corresponding cells on "metadata".
$styleSource = $dataSheet.Range("B1")
$styleSource.Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

$indexStyleSource = $dataSheet.Range("A2")
$indexStyleSource.Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("A1").Select()
